# Add a new worksheet named "Rahul" at the end of the workbook, matching
# the structure/style of the other "stepout_qc_code" sheets (Boni/Arpit/Sudhanva).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Rahul"

# Reuse the existing header style (bold, centered, bordered) from sheet "Sudhanva"
# so no new cell-style entries get created in styles.xml.
$srcSheet = $wb.Worksheets.Item("Sudhanva")
$srcSheet.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Reuse the existing date-formatted style for the "current_date" column.
$srcSheet.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("A1").Value = "team_a_name"
$ws.Range("B1").Value = "team_b_name"
$ws.Range("C1").Value = "match_id"
$ws.Range("D1").Value = "game_time"
$ws.Range("E1").Value = "current_date"
$ws.Range("F1").Value = "renumeration"

$ws.Range("A2").Value = "yuu"
$ws.Range("B2").Value = "iuu"
$ws.Range("C2").Value = 123
$ws.Range("D2").Value = 90
$ws.Range("E2").Value = 45352
$ws.Range("F2").Value = 500

# Keep the first sheet as the active one, as it was before this edit.
$wb.Worksheets.Item(1).Activate() | Out-Null
